# Auto-generated edit script: updates cryptos price/volume table
# to match the target snapshot from the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper (far outside the used range) used to push numeric-looking
# strings through Copy/PasteSpecial(xlPasteValues) so they land as
# TEXT cells (matching the source inlineStr cells) instead of being
# auto-coerced to numbers by a plain .Value assignment.
$helper = $ws.Range("Z1000")

$ws.Range("D2").Value = "57.002.15"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "2.399.20"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("E6").Value = "  +3.22%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "2.417.70"
$ws.Range("E9").Value = "  +1.83%  "
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("E11").Value = "  -1.28%  "
$ws.Range("E12").Value = "  +1.75%  "
$ws.Range("E13").Value = "  -4.45%  "
$ws.Range("D14").Value = "2.832.00"
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("D15").Value = "56.881.21"
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("E17").Value = "  +1.86%  "
$ws.Range("D18").Value = "2.362.63"
$ws.Range("E18").Value = "  -3.01%  "
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("E22").Value = "  +4.04%  "
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E28").Value = "  -3.26%  "
$ws.Range("E29").Value = "  +4.29%  "
$ws.Range("E30").Value = "  -1.69%  "
$ws.Range("D31").Value = "0.0₃0723"
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("E34").Value = "  -3.44%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("E40").Value = "  +3.35%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("E43").Value = "  +7.25%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E44").Value = "  +1.93%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("E50").Value = "  +2.23%  "
$ws.Range("E51").Value = "  +0.95%  "

# Numeric-looking text values (forced to TEXT via paste-special)
$helper.Formula = "=""505.40"""
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$helper.Formula = "=""131.84"""
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$helper.Formula = "=""0.997"""
$helper.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$helper.Formula = "=""0.554"""
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$helper.Formula = "=""0.0964"""
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$helper.Formula = "=""21.73"""
$helper.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$helper.Formula = "=""10.19"""
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$helper.Formula = "=""309.08"""
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$helper.Formula = "=""6.34"""
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$helper.Formula = "=""5.87"""
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$helper.Formula = "=""0.998"""
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$helper.Formula = "=""65.26"""
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$helper.Formula = "=""0.998"""
$helper.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$helper.Formula = "=""7.49"""
$helper.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$helper.Formula = "=""171.10"""
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$helper.Formula = "=""17.93"""
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$helper.Formula = "=""3.83"""
$helper.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$helper.Formula = "=""36.62"""
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$helper.Formula = "=""0.802"""
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$helper.Formula = "=""1.44"""
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$helper.Formula = "=""130.80"""
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$helper.Formula = "=""4.97"""
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$helper.Formula = "=""3.35"""
$helper.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$helper.Formula = "=""251.79"""
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$helper.Formula = "=""0.565"""
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$helper.Formula = "=""0.0909"""
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$helper.Formula = "=""0.0487"""
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$helper.Formula = "=""16.99"""
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4163)

# Clean up the helper cell so it leaves no trace in the sheet
$helper.ClearContents()
$excel.CutCopyMode = $false

